# Apply the "PO Forecast" sheet addition + header renames described by the diff.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" headers on the existing sheets -----
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet ------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1

# --- 3. Header row ----------------------------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the header formatting (bold, centered, bordered) from an existing sheet
$wsWeekly.Range("A1:B1").Copy() | Out-Null
$wsForecast.Range("A1:D1").PasteSpecial(-4122) | Out-Null

# --- 4. Data rows ------------------------------------------------------------
$data = @(
    @(44948.99999999999,4,-9.281505314149598,18.98929156327432),
    @(44955.99999999999,4,-9.519577964331784,19.08715695364212),
    @(44962.99999999999,5,-8.201414594473471,19.07692404391662),
    @(44969.99999999999,5,-8.363693740482033,18.89984977755833),
    @(44976.99999999999,6,-8.861460759896818,19.98724835963683),
    @(44997.99999999999,7,-6.96783951484263,21.70509387325452),
    @(45004.99999999999,8,-5.647576517877613,22.44690889020851),
    @(45011.99999999999,8,-5.150132990428413,22.54448615086459),
    @(45025.99999999999,10,-5.371479304466375,23.97371752135529),
    @(45032.99999999999,10,-4.902760307482693,24.18137106757434),
    @(45039.99999999999,11,-4.547099269441332,23.94613484454279),
    @(45046.99999999999,11,-2.898879433754258,24.55081271555313),
    @(45053.99999999999,12,-3.448719853883297,24.21216574805483),
    @(45060.99999999999,12,-1.556346898986052,26.55114618682475),
    @(45067.99999999999,13,-1.895293376355751,26.87389512925207),
    @(45074.99999999999,13,-0.09836144740700339,27.47025764468373),
    @(45081.99999999999,14,-0.1860149750072271,27.51733982234315),
    @(45088.99999999999,14,0.02681632071212453,27.96384855489089),
    @(45095.99999999999,15,-0.09922702658654373,29.03863237466594),
    @(45102.99999999999,15,0.5993553341898735,28.405236451706),
    @(45109.99999999999,16,1.173893499417531,29.40926211347719),
    @(45116.99999999999,16,1.689333820315973,29.8513842887968),
    @(45123.99999999999,17,2.254907567379863,30.85749585884521),
    @(45130.99999999999,17,3.549702671718616,30.80733029870285),
    @(45137.99999999999,18,4.174042365699745,31.68804132243287),
    @(45144.99999999999,18,4.574298114985808,33.68355634959946),
    @(45151.99999999999,19,4.884807504624115,31.6576826538704),
    @(45158.99999999999,19,5.764418886565775,32.71572985894117),
    @(45165.99999999999,20,4.927444548340753,34.27325363202704),
    @(45172.99999999999,20,6.814983506305382,33.82704079601682),
    @(45179.99999999999,21,6.485527973038133,35.69733803785198),
    @(45186.99999999999,21,7.524944557084569,34.48905815663256),
    @(45193.99999999999,22,9.039660345274546,35.85440152312398)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $rowNum = $i + 2
    $wsForecast.Cells.Item($rowNum, 1).Value = $rowValues[0]
    $wsForecast.Cells.Item($rowNum, 2).Value = $rowValues[1]
    $wsForecast.Cells.Item($rowNum, 3).Value = $rowValues[2]
    $wsForecast.Cells.Item($rowNum, 4).Value = $rowValues[3]
}

# Copy the date-cell formatting (date number format) down column A
$wsWeekly.Range("A2").Copy() | Out-Null
$wsForecast.Range("A2:A34").PasteSpecial(-4122) | Out-Null

$wsForecast.Range("A1").Select() | Out-Null
